# The sheet originally held 21 columns of test data (A:U). Columns H:T held
# a block of "internalLabel" button columns and duplicate field-label
# columns that are no longer needed; only the very last column (U, the
# "input_email" field) should remain, sliding into column H.
#
# Deleting the entire H:T column range removes that block and shifts the
# remaining column U left by 13 positions so it becomes the new column H,
# which is exactly what the target workbook looks like (dimension A1:H2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1:T2").EntireColumn.Delete()
